$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the combined format (bold font, thin box border, centered/top aligned)
# on a single cell first so only one new style record is produced, then copy
# that exact formatting onto the other cell that needs it (B1 and A2 share
# the same style; B2 stays unformatted).
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160

$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
